$wb = $excel.ActiveWorkbook

# Update the "Ngay xuat" (issue date) value shown on every sheet (shared string change)
foreach ($sheet in $wb.Worksheets) {
    $sheet.Range("B4").Value = "17/12/2025"
}

# Update refreshed financial figures on the "Note" sheet (annual 2024 total in column C,
# and Q4 2024 quarterly figure in column H) for the rows whose source data changed.
$ws = $wb.Worksheets.Item("Note")

$ws.Range("C29").Value = 683944075529
$ws.Range("H29").Value = 683944075529
$ws.Range("C35").Value = 303231966572
$ws.Range("H35").Value = 303231966572
$ws.Range("C63").Value = 7271284762
$ws.Range("H63").Value = 7271284762
$ws.Range("C64").Value = -14294204188
$ws.Range("H64").Value = -14294204188
$ws.Range("C69").Value = 5412013666
$ws.Range("H69").Value = 5412013666
$ws.Range("C71").Value = 2005061039
$ws.Range("H71").Value = 2005061039
$ws.Range("C79").Value = 864219493965
$ws.Range("H79").Value = 864219493965
$ws.Range("C80").Value = 833835440802
$ws.Range("H80").Value = 833835440802
$ws.Range("C81").Value = 30384053163
$ws.Range("H81").Value = 30384053163
$ws.Range("C90").Value = 21122532154
$ws.Range("H90").Value = 21122532154
$ws.Range("C94").Value = 15458865622
$ws.Range("H94").Value = 15458865622
$ws.Range("C105").Value = 61230124075
$ws.Range("H105").Value = 61230124075
$ws.Range("C109").Value = 187822800
$ws.Range("H109").Value = 187822800
$ws.Range("C110").Value = 83252800
$ws.Range("H110").Value = 83252800
$ws.Range("C114").Value = 52854028425
$ws.Range("H114").Value = 52854028425
$ws.Range("C115").Value = 323254050154
$ws.Range("H115").Value = 323254050154
$ws.Range("C116").Value = 302340118756
$ws.Range("H116").Value = 302340118756
$ws.Range("C119").Value = 20913931398
$ws.Range("H119").Value = 20913931398
$ws.Range("C124").Value = 1692068298858
$ws.Range("H124").Value = 530414032901
$ws.Range("C125").Value = 1685835028066
$ws.Range("H125").Value = 529887578407
$ws.Range("C126").Value = 6233270792
$ws.Range("H126").Value = 526454494
$ws.Range("C135").Value = 1425928927552
$ws.Range("H135").Value = 474997230324
$ws.Range("C136").Value = 769837275615
$ws.Range("H136").Value = 316294275282
$ws.Range("C137").Value = 647873186615
$ws.Range("H137").Value = 158257760754
$ws.Range("C138").Value = 8218465322
$ws.Range("H138").Value = 445194288
$ws.Range("C144").Value = 8758836895
$ws.Range("H144").Value = -1243997472
$ws.Range("C145").Value = 2942754831
$ws.Range("H145").Value = 1172229041
$ws.Range("C153").Value = 5816082064
$ws.Range("H153").Value = -2416226513
$ws.Range("C154").Value = 61991179349
$ws.Range("H154").Value = 14636320922
$ws.Range("C155").Value = 56620018291
$ws.Range("H155").Value = 18019420152
$ws.Range("C162").Value = 5371161058
$ws.Range("H162").Value = -3383099230
$ws.Range("C163").Value = 808131759092
$ws.Range("H163").Value = 218604262259
$ws.Range("C164").Value = 470545034765
$ws.Range("H164").Value = 132764160831
$ws.Range("C165").Value = 131021060704
$ws.Range("H165").Value = 36940505780
$ws.Range("C166").Value = 36797703957
$ws.Range("H166").Value = 9820258330
$ws.Range("C167").Value = 168936743995
$ws.Range("H167").Value = 39078995186
$ws.Range("C168").Value = 831215671
$ws.Range("H168").Value = 342132
